# Append the 09/16/2025 profit record as a new row 30 (A30:B30).
# The Date column stores the date as literal text (matching every other
# row in the sheet except the one pre-existing numeric-date outlier), so
# we force Text formatting before assigning the value to stop Excel's
# automatic date recognition from turning "09/16/2025" into a date serial.
# Resetting the style back to "Normal" afterwards keeps the new cell on
# the default (unstyled) format, just like its neighbors.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dateCell = $ws.Cells.Item(30, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "09/16/2025"
$dateCell.Style = "Normal"

$ws.Cells.Item(30, 2).Value = 15785.27
